# Fix typos in main block assignments
# The diagonal "fishN" cells (B2, C3, D4, E5, F6) were mistakenly set to
# "stimuli/bead_g.PNG" (shared string 10) instead of the intended
# "stimuli/bead_b.PNG" (shared string 12). Correct each of them.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "stimuli/bead_b.PNG"
$ws.Range("C3").Value = "stimuli/bead_b.PNG"
$ws.Range("D4").Value = "stimuli/bead_b.PNG"
$ws.Range("E5").Value = "stimuli/bead_b.PNG"
$ws.Range("F6").Value = "stimuli/bead_b.PNG"

# Update the saved view state to match: selection moved to C17:F20 and
# the sheet zoomed to 125%.
$ws.Range("C17:F20").Select()
$excel.ActiveWindow.Zoom = 125
